$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30
$ws.Range("G30").Value = 1.25
$ws.Range("S30").Value = 1.22
$ws.Range("T30").Value = 4
$ws.Range("U30").Value = 1.95
$ws.Range("V30").Value = 1.8
$ws.Range("W30").Value = 9
$ws.Range("X30").Value = 7
$ws.Range("Z30").Value = 8
$ws.Range("AG30").Value = 301
$ws.Range("AT30").Value = 4

# Row 31
$ws.Range("O31").Value = 1.36
$ws.Range("P31").Value = 3
$ws.Range("Q31").Value = 2.15
$ws.Range("R31").Value = 1.67

# Row 38
$ws.Range("H38").Value = 3.75
$ws.Range("I38").Value = 2.63
$ws.Range("L38").Value = 3.25
$ws.Range("S38").Value = 1.33
$ws.Range("T38").Value = 3.25
$ws.Range("U38").Value = 1.62
$ws.Range("V38").Value = 2.2
$ws.Range("Z38").Value = 23
$ws.Range("AB38").Value = 23
$ws.Range("AE38").Value = 13
$ws.Range("AG38").Value = 151
$ws.Range("AH38").Value = 11
$ws.Range("AI38").Value = 15
$ws.Range("AJ38").Value = 10
$ws.Range("AS38").Value = 126
$ws.Range("AT38").Value = 3.25
$ws.Range("AZ38").Value = 41

# Row 39
$ws.Range("G39").Value = 3.1
$ws.Range("H39").Value = 2.9
$ws.Range("I39").Value = 2.4
$ws.Range("L39").Value = 3.25
$ws.Range("N39").Value = 7.5
$ws.Range("O39").Value = 1.36
$ws.Range("P39").Value = 3
$ws.Range("Q39").Value = 2.25
$ws.Range("R39").Value = 1.62
$ws.Range("U39").Value = 1.83
$ws.Range("V39").Value = 1.83
$ws.Range("W39").Value = 8.5
$ws.Range("AC39").Value = 7.5
$ws.Range("AD39").Value = 6
$ws.Range("AG39").Value = 301
$ws.Range("AJ39").Value = 10
$ws.Range("AK39").Value = 23
$ws.Range("AL39").Value = 21
$ws.Range("AM39").Value = 34
$ws.Range("AS39").Value = 201
$ws.Range("AY39").Value = 26
$ws.Range("BB39").Value = 201

# Row 46
$ws.Range("G46").Value = 3.25
$ws.Range("H46").Value = 3.9
$ws.Range("I46").Value = 2.05
$ws.Range("J46").Value = 3.6
$ws.Range("L46").Value = 2.6
$ws.Range("N46").Value = 19
$ws.Range("O46").Value = 1.17
$ws.Range("S46").Value = 1.29
$ws.Range("T46").Value = 3.5
$ws.Range("U46").Value = 1.5
$ws.Range("W46").Value = 13
$ws.Range("X46").Value = 19
$ws.Range("Y46").Value = 11
$ws.Range("Z46").Value = 34
$ws.Range("AB46").Value = 23
$ws.Range("AL46").Value = 15
$ws.Range("AM46").Value = 21
$ws.Range("AN46").Value = 5.5
$ws.Range("AT46").Value = 3.5
$ws.Range("AW46").Value = 4.5
$ws.Range("AZ46").Value = 34

$wb.Save()
